# User Story Specs.xlsx - update "example:" story block from the old bike-rack
# placeholder text to the real GetRichTwitch Twitch-bot user stories, and add
# a brand-new user-story row (row 10) for the "!stats" command.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: "!slap" user story (replaces the old bike-rack example text) ---
$ws.Range("A8").Value = 'I, |slartibartfast|, am a Twitch user that wants to have fun in Twitch chat. I type "!slap" and then the name of another user in chat to slap them.'
$ws.Range("B8").Value = 'The bot slaps the wanted user, typing @[username] You Have Been Slapped! '
$ws.Range("C8").Value = 'Searches the running list of usernames that have posted in the chat, tells slaper whether slapee is in chat.'
$ws.Range("D8").Value = 'keep a tally of who has been slapped the most'
$ws.Range("E8").Value = 'An actual hand coming out of the computer to slap someone, because it is all digital and fake internet points.'
$ws.Range("F8").Value = 'Levi'
$ws.Range("F8").WrapText = $true

$ws.Rows.Item(8).RowHeight = 106.5

# --- Row 9: second line of the "!slap" story (Could Have / Won't Have cells) ---
$ws.Range("C9").ClearContents()
$ws.Range("D9").Value = 'permissions, only certain users can slap'

$ws.Rows.Item(9).RowHeight = 41.25

# --- Row 10: new "!stats" user story ---
$ws.Range("A10").Value = 'I, XxMaestroChefxX is a Twitch user that wants to see their status. They type "!stats" in chat.'
$ws.Range("B10").Value = 'The bot responds with information on their status including name, if they are a mod, and badges.'
$ws.Range("D10").Value = 'also includes information on their points, coins, and other chat games information.'
$ws.Range("C10").Value = 'stores information in a file for data mining/management for use of account using bot.'
$ws.Range("E10").Value = 'Personal information will not be stored.'
$ws.Range("F10").Value = 'Levi'
$ws.Range("F10").WrapText = $true

$ws.Rows.Item(10).RowHeight = 131.25

# A10 gets its own wrap/vertical-center style (no horizontal alignment, unlike A8/A9).
$ws.Range("A10").WrapText = $true
$ws.Range("A10").VerticalAlignment = -4108

# --- Selection moves to E10 (last cell touched) ---
$ws.Range("E10").Select()
